# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "581.48") are not coerced into numbers, then restore the
    # default "Normal" style so no stray formatting is introduced.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '61.810.48'
Set-TextValue $ws.Range('E2') '  +1.31%  '
Set-TextValue $ws.Range('D3') '3.460.45'
Set-TextValue $ws.Range('E3') '  +2.83%  '
Set-TextValue $ws.Range('E4') '  -0.01%  '
Set-TextValue $ws.Range('D5') '581.48'
Set-TextValue $ws.Range('E5') '  +1.86%  '
Set-TextValue $ws.Range('D6') '148.79'
Set-TextValue $ws.Range('E6') '  +9.57%  '
Set-TextValue $ws.Range('D7') '3.463.43'
Set-TextValue $ws.Range('E7') '  +3.04%  '
Set-TextValue $ws.Range('E8') '  +0.03%  '
Set-TextValue $ws.Range('D9') '0.474'
Set-TextValue $ws.Range('E9') '  +1.14%  '
Set-TextValue $ws.Range('D10') '7.81'
Set-TextValue $ws.Range('E10') '  +4.19%  '
Set-TextValue $ws.Range('E11') '  +1.77%  '
Set-TextValue $ws.Range('D12') '0.391'
Set-TextValue $ws.Range('E12') '  +1.78%  '
Set-TextValue $ws.Range('D13') '4.054.70'
Set-TextValue $ws.Range('E13') '  +2.90%  '
Set-TextValue $ws.Range('D14') '28.18'
Set-TextValue $ws.Range('E14') '  +8.46%  '
Set-TextValue $ws.Range('E15') '  -0.46%  '
Set-TextValue $ws.Range('D16') '0.0000175'
Set-TextValue $ws.Range('E16') '  +1.93%  '
Set-TextValue $ws.Range('D17') '3.458.47'
Set-TextValue $ws.Range('E17') '  +2.78%  '
Set-TextValue $ws.Range('D18') '61.865.62'
Set-TextValue $ws.Range('E18') '  +1.19%  '
Set-TextValue $ws.Range('D19') '6.34'
Set-TextValue $ws.Range('E19') '  +9.06%  '
Set-TextValue $ws.Range('D20') '14.37'
Set-TextValue $ws.Range('E20') '  +2.92%  '
Set-TextValue $ws.Range('E21') '  +2.61%  '
Set-TextValue $ws.Range('D22') '385.93'
Set-TextValue $ws.Range('E22') '  +2.40%  '
Set-TextValue $ws.Range('D23') '0.569'
Set-TextValue $ws.Range('E23') '  +2.88%  '
Set-TextValue $ws.Range('D24') '3.600.44'
Set-TextValue $ws.Range('E24') '  +2.89%  '
Set-TextValue $ws.Range('D25') '72.72'
Set-TextValue $ws.Range('E25') '  +2.38%  '
Set-TextValue $ws.Range('E26') '  +1.08%  '
Set-TextValue $ws.Range('D27') '0.999'
Set-TextValue $ws.Range('E27') '  -0.08%  '
Set-TextValue $ws.Range('E28') '  -1.80%  '
Set-TextValue $ws.Range('D29') '0.181'
Set-TextValue $ws.Range('E29') '  +9.44%  '
Set-TextValue $ws.Range('D30') '7.81'
Set-TextValue $ws.Range('E30') '  +4.15%  '
Set-TextValue $ws.Range('E31') '  -0.33%  '
Set-TextValue $ws.Range('D32') '1.52'
Set-TextValue $ws.Range('E32') '  -13.58%  '
Set-TextValue $ws.Range('D33') '8.25'
Set-TextValue $ws.Range('E33') '  +1.26%  '
Set-TextValue $ws.Range('E34') '  +1.82%  '
Set-TextValue $ws.Range('D36') '24.00'
Set-TextValue $ws.Range('E36') '  +1.76%  '
Set-TextValue $ws.Range('E37') '  +4.59%  '
Set-TextValue $ws.Range('E38') '  +0.64%  '
Set-TextValue $ws.Range('E39') '  +2.53%  '
Set-TextValue $ws.Range('D40') '166.09'
Set-TextValue $ws.Range('E40') '  +0.66%  '
Set-TextValue $ws.Range('D41') '0.0791'
Set-TextValue $ws.Range('E41') '  +5.33%  '
Set-TextValue $ws.Range('D42') '26.24'
Set-TextValue $ws.Range('E42') '  +9.92%  '
Set-TextValue $ws.Range('E44') '  +0.00%  '
Set-TextValue $ws.Range('E45') '  +2.20%  '
Set-TextValue $ws.Range('B46') 'Filecoin'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D46') '4.49'
Set-TextValue $ws.Range('E46') '  +2.19%  '
Set-TextValue $ws.Range('B47') 'Stacks'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D47') '1.73'
Set-TextValue $ws.Range('E47') '  +1.42%  '
Set-TextValue $ws.Range('D48') '1.17'
Set-TextValue $ws.Range('E48') '  -1.85%  '
Set-TextValue $ws.Range('D49') '2.601.49'
Set-TextValue $ws.Range('E49') '  +10.55%  '
Set-TextValue $ws.Range('D50') '6.98'
Set-TextValue $ws.Range('E50') '  +2.70%  '
Set-TextValue $ws.Range('D51') '23.33'
Set-TextValue $ws.Range('E51') '  +0.59%  '
